$wb = $excel.ActiveWorkbook

# Copy the "Portugal" sheet to create the new "Slovakia" sheet, placing it
# right after "Portugal" (at the end of the tab strip).
$portugal = $wb.Worksheets.Item("Portugal")
$portugal.Copy($null, $portugal)

$ws = $wb.Worksheets.Item($portugal.Index + 1)
$ws.Name = "Slovakia"

# Update the market-specific content on the new sheet.
$ws.Range("B2").Value = "Slovakia Market"

# Reset B4 back to the default (unstyled) cell before typing the new value,
# so the pasted-from-Portugal border/fill formatting is dropped - matching
# how the value was entered by hand on the new sheet.
$ws.Range("B4").Style = "Normal"
$ws.Range("B4").Value = "NGC-2930/T3180"

# Leave the Portugal sheet with everything selected, then make the new
# Slovakia sheet the active tab with its last populated cell selected.
$portugal.Cells.Select()

$ws.Activate()
$ws.Range("A12").Select()
